# Insert a new data row at row 43 (shifts existing rows 43:99 down to 44:100)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 43:99 down to 44:100 by inserting a new row at position 43.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record.
# Most fields mirror the row immediately below it (the former row 43,
# now shifted to row 44) except the date and the price columns.
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value = (Get-Date -Year 2022 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(43, 5).Value = 15
$ws.Cells.Item(43, 6).Value = 100112038
$ws.Cells.Item(43, 7).Value = "Cebollín baby"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 300
$ws.Cells.Item(43, 11).Value = 2400
$ws.Cells.Item(43, 12).Value = 2500
$ws.Cells.Item(43, 13).Value = 2450
$ws.Cells.Item(43, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 1225
$ws.Cells.Item(43, 17).Value = 2
$ws.Cells.Item(43, 18).Value = "Hortaliza"

# Copy the date cell's number format from the row below (style "s=2" -> custom
# date/time format) so the new row matches the rest of the date column.
$ws.Cells.Item(44, 4).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4122) | Out-Null
